$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the text of the 4th-iteration task description (A31)
$ws.Range("A31").Value = "4. iterace - další implementace (vybrat stůl, upravit rezervaci)"

# 2. Bump the hours logged for that task from 1 to 2 (B4 total recalculates automatically)
$ws.Range("B31").Value = 2

# 3. Give A32 the same (blank) formatting as the task column above it (A18:A31use style index 13)
#    by inserting a row above it (inherits formatting from row 31) and then removing the
#    now-duplicated blank row that got pushed down, so the row count/content stays identical.
$ws.Rows("32").Insert(-4121)
$ws.Rows("33").Delete()

# 4. Move the active selection to F19
$ws.Range("F19").Select()
